$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 15 (TI_EQUIPES_TOURNOIS): update "Auteur" then "But" columns
# (order matters for shared-string index allocation)
$ws.Range("D15").Value = "Jonathan"
$ws.Range("B15").Value = "*ID - Insert transaction"

# Row 18 (TI_INSCRIPTIONS): update "Auteur" column
$ws.Range("D18").Value = "Jonathan"

# Update the view state: scroll position and active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
$ws.Range("D18").Select()
